$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 167
$ws.Range("D2").Value = 31
$ws.Range("E2").Value = 136

$ws.Range("C3").Value = 172
$ws.Range("D3").Value = 35

$ws.Range("C4").Value = 196
$ws.Range("D4").Value = 49

$ws.Range("C5").Value = 187
$ws.Range("D5").Value = 39
